# Update New Orleans workbook:
#  1. Add a "State" column to hotel_info (between Hotel_Name and City) with value "Louisiana".
#  2. Reorder the sheet tabs so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "State" column into hotel_info ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns("C").Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Move review_info so it becomes the first sheet ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($wb.Worksheets.Item(1))
